$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Price cells that receive a new value to Text format first so that
# strings such as "305.64" or "49.00" are kept verbatim (not coerced into numbers),
# matching how these values are stored as text in the original workbook.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.808.86'
$ws.Range("E2").Value = '  -3.04%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.616.57'
$ws.Range("E3").Value = '  -3.47%  '

$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("E5").Value = '  +0.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '305.64'
$ws.Range("E6").Value = '  -2.66%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3890'
$ws.Range("E7").Value = '  -0.25%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3825'
$ws.Range("E8").Value = '  -2.84%  '

$ws.Range("E9").Value = '  +0.03%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = '1.351'
$ws.Range("E10").Value = '  -2.97%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").Value = '49.00'
$ws.Range("E11").Value = '  -5.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08388'
$ws.Range("E12").Value = '  -2.76%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.59'
$ws.Range("E13").Value = '  -5.78%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.938'
$ws.Range("E14").Value = '  -4.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001267'
$ws.Range("E15").Value = '  -3.48%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.391'
$ws.Range("E16").Value = '  -4.25%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.618.58'
$ws.Range("E17").Value = '  -4.22%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.63'
$ws.Range("E18").Value = '  -1.10%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06903'
$ws.Range("E19").Value = '  -2.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.75'
$ws.Range("E20").Value = '  -3.60%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.793'
$ws.Range("E21").Value = '  -3.61%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.005'
$ws.Range("E22").Value = '  -0.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.30'
$ws.Range("E23").Value = '  -4.99%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.838.21'
$ws.Range("E24").Value = '  -2.92%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.413'
$ws.Range("E25").Value = '  +1.46%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.792'
$ws.Range("E26").Value = '  +2.68%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.98'
$ws.Range("E27").Value = '  -5.57%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.40'
$ws.Range("E28").Value = '  -3.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '138.95'
$ws.Range("E29").Value = '  -5.42%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.242'
$ws.Range("E30").Value = '  -9.10%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.717'
$ws.Range("E31").Value = '  -3.96%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.479'
$ws.Range("E32").Value = '  -1.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.799.58'
$ws.Range("E33").Value = '  -3.99%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07943'
$ws.Range("E34").Value = '  -4.94%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02867'
$ws.Range("E35").Value = '  -5.49%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9566'
$ws.Range("E36").Value = '  -1.36%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.545'
$ws.Range("E37").Value = '  -3.61%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2637'
$ws.Range("E38").Value = '  -6.13%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09094'
$ws.Range("E39").Value = '  -4.24%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.33'
$ws.Range("E40").Value = '  +0.17%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.47'
$ws.Range("E41").Value = '  -0.33%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.423'
$ws.Range("E42").Value = '  -7.78%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7407'
$ws.Range("E43").Value = '  -5.86%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.76'
$ws.Range("E44").Value = '  -4.26%  '

$ws.Range("E45").Value = '  -4.24%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.421'
$ws.Range("E46").Value = '  -5.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.047'
$ws.Range("E47").Value = '  -3.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.004'
$ws.Range("E48").Value = '  +0.12%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08221'
$ws.Range("E49").Value = '  -4.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.87'
$ws.Range("E50").Value = '  -3.81%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.242'
$ws.Range("E51").Value = '  -5.61%  '

